$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column E width
$ws.Columns.Item(5).ColumnWidth = 17.43

# Row heights for rows 1-5
$ws.Rows.Item(1).RowHeight = 16.5
$ws.Rows.Item(2).RowHeight = 30.75
$ws.Rows.Item(3).RowHeight = 15.75
$ws.Rows.Item(4).RowHeight = 15.75
$ws.Rows.Item(5).RowHeight = 30.75

# New column E values (order matters for shared string table indices)
$ws.Range("E23").Value = "Initial"
$ws.Range("E2").Value = "Chart initial"
$ws.Range("E9").Value = "Blood test, initial"
$ws.Range("E14").Value = "Blood chemistry, initial"
$ws.Range("E35").Value = "Initial, diet"
$ws.Range("E36").Value = "Intial, urine"
$ws.Range("E39").Value = "Initial"
$ws.Range("E42").Value = "Initial, blood chemistry"
$ws.Range("E49").Value = "Lowering Sodium Intake"
$ws.Range("E52").Value = "Blocking AngII, 70%"
$ws.Range("E55").Value = "Lowering Sodium and Blocking AngII"

# View changes
$ws.Application.ActiveWindow.ScrollRow = 25
$ws.Range("G49").Select()
